{"js": "// Resume content tweaks: soften a handful of bullet/skill phrasings.\n// Each entry is [oldText, newText] where oldText is the exact, unique\n// text of a single run in the document.\nconst replacements = [\n  [\n    \"git, jest, vitest, React Testing Library, TurboRepo, CI/CD Pipelines and Actions, Redux / Redux Saga\",\n    \"git, jest, vitest, React Testing Library, TurboRepo, CI/CD Pipelines, Redux / Redux Saga\"\n  ],\n  [\n    \"Serves as the primary point of oversight and escalation for numerous projects\",\n    \"Serves as the primary point of oversight and escalation for several projects\"\n  ],\n  [\n    \"Spearheads the development and architecture of several internal company projects\",\n    \"Spearheads the development and architecture of multiple internal company projects\"\n  ],\n  [\n    \"Developed a large-scale digital health application in AngularJS, providing services to healthcare consumers and providers\",\n    \"Contributed to the development of a large-scale digital health application in AngularJS, providing services to healthcare consumers and providers\"\n  ],\n  [\n    \"Contributed to the development of an internal web application designed to streamline weekly status reporting for developers\",\n    \"Played a key role in leading and building an internal web application designed to streamline weekly status reporting for developers\"\n  ],\n  [\n    \"Member of the team responsible for developing an internal software system utilized by the sales and dispatch team, which tracks and aligns customer loads with carrier trucks\",\n    \"Member of the development team responsible for building an internal software system designed to track and align customer loads with carrier trucks, utilized by the sales and dispatch unit\"\n  ],\n  [\n    \"Primarily responsible for most user interface development, collaborating closely with the design team to bring their concepts to life\",\n    \"Primarily responsible for most user interface development, collaborating closely with the design team to bring concepts to life\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Resume content tweaks: soften a handful of bullet/skill phrasings.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"git, jest, vitest, React Testing Library, TurboRepo, CI/CD Pipelines and Actions, Redux / Redux Saga\"\n        New = \"git, jest, vitest, React Testing Library, TurboRepo, CI/CD Pipelines, Redux / Redux Saga\"\n    },\n    @{\n        Old = \"Serves as the primary point of oversight and escalation for numerous projects\"\n        New = \"Serves as the primary point of oversight and escalation for several projects\"\n    },\n    @{\n        Old = \"Spearheads the development and architecture of several internal company projects\"\n        New = \"Spearheads the development and architecture of multiple internal company projects\"\n    },\n    @{\n        Old = \"Developed a large-scale digital health application in AngularJS, providing services to healthcare consumers and providers\"\n        New = \"Contributed to the development of a large-scale digital health application in AngularJS, providing services to healthcare consumers and providers\"\n    },\n    @{\n        Old = \"Contributed to the development of an internal web application designed to streamline weekly status reporting for developers\"\n        New = \"Played a key role in leading and building an internal web application designed to streamline weekly status reporting for developers\"\n    },\n    @{\n        Old = \"Member of the team responsible for developing an internal software system utilized by the sales and dispatch team, which tracks and aligns customer loads with carrier trucks\"\n        New = \"Member of the development team responsible for building an internal software system designed to track and align customer loads with carrier trucks, utilized by the sales and dispatch unit\"\n    },\n    @{\n        Old = \"Primarily responsible for most user interface development, collaborating closely with the design team to bring their concepts to life\"\n        New = \"Primarily responsible for most user interface development, collaborating closely with the design team to bring concepts to life\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    # 0=Wrap:wdFindStop/wdFindContinue:1/wdFindAsk:2 use 1; 2=Replace:wdReplaceAll\n    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
